$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $oldStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $oldStyle
}

$ws.Range('D2').Value = '76.315.81'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.033.60'
$ws.Range('E3').Value = '  +3.33%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws.Range('D5') '200.01'
$ws.Range('E5').Value = '  -1.33%  '
Set-TextValue $ws.Range('D6') '623.98'
$ws.Range('E6').Value = '  +4.24%  '
Set-TextValue $ws.Range('D7') '1.00'
Set-TextValue $ws.Range('D8') '0.549'
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('E9').Value = '  +2.93%  '
$ws.Range('D10').Value = '3.032.41'
$ws.Range('E10').Value = '  +3.32%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  -0.99%  '
Set-TextValue $ws.Range('D13') '5.25'
$ws.Range('D14').Value = '3.591.45'
$ws.Range('E14').Value = '  +3.40%  '
Set-TextValue $ws.Range('D15') '29.05'
$ws.Range('E15').Value = '  +2.89%  '
$ws.Range('D16').Value = '76.289.20'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '3.030.29'
$ws.Range('E18').Value = '  +3.33%  '
Set-TextValue $ws.Range('D19') '13.47'
$ws.Range('E19').Value = '  +1.42%  '
Set-TextValue $ws.Range('D20') '8.95'
$ws.Range('E20').Value = '  -0.25%  '
Set-TextValue $ws.Range('D21') '373.94'
$ws.Range('E21').Value = '  -0.11%  '
Set-TextValue $ws.Range('D22') '2.29'
$ws.Range('E22').Value = '  -0.71%  '
Set-TextValue $ws.Range('D23') '4.34'
$ws.Range('E23').Value = '  +0.52%  '
Set-TextValue $ws.Range('D24') '73.07'
$ws.Range('E24').Value = '  +1.62%  '
$ws.Range('D25').Value = '3.191.19'
$ws.Range('E25').Value = '  +3.31%  '
Set-TextValue $ws.Range('D26') '0.998'
$ws.Range('E26').Value = '  -0.45%  '
Set-TextValue $ws.Range('D27') '4.36'
$ws.Range('E27').Value = '  +1.14%  '
Set-TextValue $ws.Range('D28') '9.78'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('E30').Value = '  +0.22%  '
Set-TextValue $ws.Range('D31') '8.22'
$ws.Range('E31').Value = '  +4.78%  '
Set-TextValue $ws.Range('D32') '1.40'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  +5.38%  '
Set-TextValue $ws.Range('D34') '493.29'
$ws.Range('E34').Value = '  -2.07%  '
Set-TextValue $ws.Range('D35') '1.00'
$ws.Range('E35').Value = '  -0.04%  '
Set-TextValue $ws.Range('D36') '20.54'
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('E37').Value = '  -1.41%  '
$ws.Range('E38').Value = '  +2.04%  '
Set-TextValue $ws.Range('D39') '0.383'
$ws.Range('E39').Value = '  +0.24%  '
Set-TextValue $ws.Range('D40') '0.115'
$ws.Range('E40').Value = '  +0.52%  '
Set-TextValue $ws.Range('D41') '188.93'
$ws.Range('E41').Value = '  +3.24%  '
Set-TextValue $ws.Range('D42') '0.106'
$ws.Range('E42').Value = '  -3.63%  '
$ws.Range('E43').Value = '  -0.02%  '
Set-TextValue $ws.Range('D44') '0.796'
$ws.Range('E44').Value = '  +20.39%  '
Set-TextValue $ws.Range('D45') '5.09'
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D46') '41.99'
$ws.Range('E46').Value = '  +4.41%  '
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D47') '1.26'
$ws.Range('E47').Value = '  +4.41%  '
$ws.Range('E48').Value = '  -1.43%  '
Set-TextValue $ws.Range('D49') '2.45'
$ws.Range('E49').Value = '  +4.00%  '
Set-TextValue $ws.Range('D50') '0.604'
$ws.Range('E50').Value = '  +3.66%  '
Set-TextValue $ws.Range('D51') '3.89'
$ws.Range('E51').Value = '  +3.34%  '
